# Generate Report for Handoff
# Update file id GUID references, hash, and timestamps across the three sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "ea7cd848-33b1-44ac-9134-4e2181dcbd27"
$newGuid = "fed28c1a-55ce-41a9-9fd9-b5c3ed70f9f0"
$oldHash = "062e217eed3d7f716793398f59d71af946de1ce9"
$newHash = "271fe60cb78713f1a0a8bdfde164f2e7f81cd03d"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 14:58:36"

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-16 14:58:31"

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
